$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-03-01 06:48:19"
$ws.Range("O2").Value = "-1.1 °C"
$ws.Range("E3").Value = "2026-03-01 06:48:22"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "97%"
$ws.Range("N3").Value = "-4.2 °C 6:04 TU"
$ws.Range("E4").Value = "2026-03-01 06:48:24"
$ws.Range("E5").Value = "2026-03-01 06:48:26"
$ws.Range("K5").Value = "-0.1 MJ/m2"
$ws.Range("O5").Value = "-3.9 °C"
$ws.Range("E6").Value = "2026-03-01 06:48:29"
$ws.Range("E7").Value = "2026-03-01 06:48:31"
$ws.Range("E8").Value = "2026-03-01 06:48:34"
$ws.Range("E9").Value = "2026-03-01 06:48:36"
$ws.Range("N9").Value = "11.1 °C 6:29 TU"
$ws.Range("E10").Value = "2026-03-01 06:48:39"
$ws.Range("N10").Value = "3.5 °C 6:27 TU"
$ws.Range("O10").Value = "6.0 °C"
$ws.Range("E11").Value = "2026-03-01 06:48:41"
$ws.Range("N11").Value = "5.8 °C 6:00 TU"
$ws.Range("O11").Value = "6.2 °C"
$ws.Range("E12").Value = "2026-03-01 06:48:43"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "76%"
$ws.Range("N12").Value = "8.5 °C 6:16 TU"
$ws.Range("O12").Value = "10.2 °C"
$ws.Range("E13").Value = "2026-03-01 06:48:46"
$ws.Range("J13").Value = "1026.3 hPa"
$ws.Range("N13").Value = "3.2 °C 6:29 TU"
$ws.Range("O13").Value = "4.2 °C"
$ws.Range("E14").Value = "2026-03-01 06:48:48"
$ws.Range("O14").Value = "10.7 °C"
$ws.Range("E15").Value = "2026-03-01 06:48:50"
$ws.Range("O15").Value = "8.6 °C"
$ws.Range("E16").Value = "2026-03-01 06:48:53"
$ws.Range("N16").Value = "-7.2 °C 6:06 TU"
$ws.Range("O16").Value = "-5.3 °C"
$ws.Range("E17").Value = "2026-03-01 06:48:55"
$ws.Range("I17").Value = "0.1 mm"
$ws.Range("E18").Value = "2026-03-01 06:48:58"
$ws.Range("E19").Value = "2026-03-01 06:49:00"
$ws.Range("N19").Value = "5.9 °C 6:10 TU"
$ws.Range("O19").Value = "6.0 °C"
$ws.Range("E20").Value = "2026-03-01 06:49:02"
$ws.Range("L20").Value = "11.5 km/h - 219º 6:12 TU"
$ws.Range("N20").Value = "-5.0 °C 6:25 TU"
$ws.Range("O20").Value = "-3.2 °C"
$ws.Range("E21").Value = "2026-03-01 06:49:05"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "89%"
$ws.Range("N21").Value = "4.9 °C 6:25 TU"
$ws.Range("O21").Value = "6.3 °C"
$ws.Range("E22").Value = "2026-03-01 06:49:07"
$ws.Range("N22").Value = "-6.6 °C 6:19 TU"
$ws.Range("E23").Value = "2026-03-01 06:49:10"
$ws.Range("N23").Value = "-4.5 °C 6:16 TU"
$ws.Range("E24").Value = "2026-03-01 06:49:12"
$ws.Range("M24").Value = "6.3 °C 6:29 TU"
$ws.Range("O24").Value = "4.7 °C"
$ws.Range("E25").Value = "2026-03-01 06:49:15"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "93%"
$ws.Range("N25").Value = "-3.4 °C 6:29 TU"
$ws.Range("O25").Value = "-2.3 °C"
$ws.Range("E26").Value = "2026-03-01 06:49:17"
$ws.Range("E27").Value = "2026-03-01 06:49:19"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "96%"
$ws.Range("N27").Value = "-2.6 °C 6:29 TU"
$ws.Range("O27").Value = "-1.6 °C"
$ws.Range("E28").Value = "2026-03-01 06:49:22"
$ws.Range("N28").Value = "8.3 °C 6:29 TU"
$ws.Range("E29").Value = "2026-03-01 06:49:24"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "81%"
$ws.Range("N29").Value = "6.6 °C 6:29 TU"
$ws.Range("O29").Value = "9.0 °C"
$ws.Range("E30").Value = "2026-03-01 06:49:26"
$ws.Range("O30").Value = "10.1 °C"
$ws.Range("E31").Value = "2026-03-01 06:49:29"
$ws.Range("E32").Value = "2026-03-01 06:49:31"
$ws.Range("M32").Value = "5.1 °C 6:10 TU"
$ws.Range("O32").Value = "2.9 °C"
$ws.Range("E33").Value = "2026-03-01 06:49:34"
$ws.Range("N33").Value = "3.6 °C 6:07 TU"
$ws.Range("E34").Value = "2026-03-01 06:49:36"
$ws.Range("N34").Value = "-0.8 °C 6:19 TU"
$ws.Range("E35").Value = "2026-03-01 06:49:39"
$ws.Range("E36").Value = "2026-03-01 06:49:41"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "74%"
$ws.Range("M36").Value = "12.7 °C 6:01 TU"
$ws.Range("O36").Value = "10.3 °C"
$ws.Range("E37").Value = "2026-03-01 06:49:44"
$ws.Range("I37").Value = "0.4 mm"
$ws.Range("N37").Value = "6.1 °C 6:04 TU"
$ws.Range("E38").Value = "2026-03-01 06:49:46"
$ws.Range("E39").Value = "2026-03-01 06:49:48"
$ws.Range("E40").Value = "2026-03-01 06:49:50"
$ws.Range("J40").Value = "1025.8 hPa"
$ws.Range("O40").Value = "6.8 °C"
$ws.Range("E41").Value = "2026-03-01 06:49:53"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "88%"
$ws.Range("N41").Value = "10.6 °C 6:12 TU"
$ws.Range("O41").Value = "11.6 °C"
$ws.Range("E42").Value = "2026-03-01 06:49:55"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H42").Value = "88%"
$ws.Range("O42").Value = "8.1 °C"
$ws.Range("E43").Value = "2026-03-01 06:49:57"
$ws.Range("E44").Value = "2026-03-01 06:50:00"
$ws.Range("N44").Value = "-4.3 °C 6:05 TU"
$ws.Range("O44").Value = "-2.9 °C"
$ws.Range("E45").Value = "2026-03-01 06:50:02"
$ws.Range("H45").NumberFormat = "@"
$ws.Range("H45").Value = "94%"
$ws.Range("J45").Value = "1027.1 hPa"
$ws.Range("N45").Value = "1.2 °C 6:24 TU"
$ws.Range("O45").Value = "3.3 °C"
$ws.Range("E46").Value = "2026-03-01 06:50:04"
$ws.Range("O46").Value = "7.8 °C"
